# journey testing - changes for spanish cert filepaths, added sep cert user input
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# The old "SpanishHeaders" section header (row 51) is replaced by a new "Regex" section header
$ws.Range("A51").Value = "Regex"

# A new "SeparateProducts" regex entry is added under the new "Regex" section (row 52)
$ws.Range("A52").Value = "SeparateProducts"
$ws.Range("B52").Value = "{\Wproduct\W:\W(\d+\s*)+\W}"

# Reflect the scrolled / selected state left behind after the edit
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B61").Select()
